# The commit renames the "Values" column header to "Value" and — because
# the table's underlying column order had drifted out of sync with the
# worksheet's left-to-right column order — re-creates the table so its
# column ids line up with the sheet columns again (id 1 -> col A "Row",
# id 2 -> col B "Value"). It also bumps the table style by one
# (TableStyleLight8 -> TableStyleLight9) and moves the active selection
# from A5 to F6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing table definition (keeps the cell data/formatting).
$lo = $ws.ListObjects.Item(1)
$lo.Unlist()

# Rename the header cell: "Values" -> "Value" (the data rows are untouched).
$ws.Range("B1").Value = "Value"

# Re-create the table over the same range so the column ids realign with
# the sheet's physical column order, then restore its name/style.
$newlo = $ws.ListObjects.Add(1, $ws.Range("A1:B4"), $null, 1)
$newlo.Name = "Table1"
$newlo.TableStyle = "TableStyleLight9"

# Move the active cell selection from A5 to F6.
$ws.Range("F6").Select()
